$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.797.72'
$ws.Range("E2").Value = '  -0.47%  '

$ws.Range("D3").Value = '3.800.19'
$ws.Range("E3").Value = '  -1.45%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.09%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '704.74'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.35%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '170.71'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.46%  '

$ws.Range("D7").Value = '3.797.92'
$ws.Range("E7").Value = '  -1.45%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("E9").Value = '  -1.11%  '

$ws.Range("E10").Value = '  -1.77%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '7.39'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.53%  '

$ws.Range("E12").Value = '  -1.23%  '

$ws.Range("E13").Value = '  -2.05%  '

$ws.Range("E14").Value = '  -1.29%  '

$ws.Range("D15").Value = '4.436.27'
$ws.Range("E15").Value = '  -1.62%  '

$ws.Range("D16").Value = '3.797.59'
$ws.Range("E16").Value = '  -1.68%  '

$ws.Range("D17").Value = '70.733.61'
$ws.Range("E17").Value = '  -0.66%  '

$ws.Range("E19").Value = '  -1.95%  '

$ws.Range("E20").Value = '  -2.07%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '495.48'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.08%  '

$ws.Range("E22").Value = '  -4.40%  '

$ws.Range("E23").Value = '  +0.80%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '85.10'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.33%  '

$ws.Range("E25").Value = '  -1.08%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '12.10'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -1.95%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '10.46'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -1.19%  '

$ws.Range("D28").Value = '3.947.96'
$ws.Range("E28").Value = '  -1.71%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.03%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '2.05'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -4.64%  '

$ws.Range("E31").Value = '  -2.77%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '7.34'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -3.77%  '

$ws.Range("E33").Value = '  -3.95%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '29.12'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -2.20%  '

$ws.Range("E35").Value = '  -2.57%  '

$ws.Range("B36").Value = 'Binance-PegBSC-USD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -0.04%  '

$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '9.09'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -2.40%  '

$ws.Range("B38").Value = 'RenzoRestakedETH'
$ws.Range("C38").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D38").Value = '3.767.63'
$ws.Range("E38").Value = '  -1.08%  '

$ws.Range("E39").Value = '  -3.51%  '

$ws.Range("E40").Value = '  +1.80%  '

$ws.Range("E41").Value = '  -3.04%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '5.94'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -2.14%  '

$ws.Range("E43").Value = '  -3.72%  '

$ws.Range("E44").Value = '  -0.03%  '

$ws.Range("E45").Value = '  +0.02%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.000320'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +4.54%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '164.62'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.04%  '

$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '426.38'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +1.37%  '

$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '48.80'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.19%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '8.60'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.96%  '

$ws.Range("E51").Value = '  -1.75%  '
